$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.011191964149475
$ws.Range("B1").Value = 2.127410411834717
$ws.Range("C1").Value = 5.880775451660156
$ws.Range("D1").Value = 1.105409979820251
$ws.Range("E1").Value = 1.156762361526489
